# Swap the match-data (columns B..AD) between pairs of rows whose
# home/away team order had been mixed up in the source feed.
# Column A (the row's sequential id / index) is intentionally left in
# place; only the "payload" columns (id, Div..PL_AhUnder) are exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First and last payload columns to swap (inclusive).
$firstCol = 2   # B
$lastCol  = 30  # AD

# Row-number pairs that need their contents exchanged.
$pairs = @(
    @(17, 18),
    @(41, 42),
    @(78, 79),
    @(89, 90),
    @(91, 92),
    @(103, 104),
    @(108, 109),
    @(135, 136),
    @(231, 232)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Snapshot every cell's value (and whether it actually holds data)
    # from both rows before writing anything back, since columns will be
    # overwritten in place as we iterate.
    $row1Vals = @{}
    $row2Vals = @{}
    $row1Has  = @{}
    $row2Has  = @{}

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)

        $row1Vals[$col] = $c1.Value2
        $row2Vals[$col] = $c2.Value2
        $row1Has[$col]  = -not $c1.Value2.Equals($null)
        $row2Has[$col]  = -not $c2.Value2.Equals($null)
    }

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)

        if ($row2Has[$col]) {
            $c1.Value = $row2Vals[$col]
        } else {
            $c1.ClearContents()
        }

        if ($row1Has[$col]) {
            $c2.Value = $row1Vals[$col]
        } else {
            $c2.ClearContents()
        }
    }
}
